# Updates cryptocurrency Price (D) and Volume(1h) (E) columns to match
# the latest scrape. Numeric-looking Price values are written with a
# leading apostrophe so Excel keeps storing them as text (matching the
# original inline-string cells) instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.702.92'
$ws.Range('E2').Value = '  -3.06%  '
# Row 3
$ws.Range('D3').Value = '3.166.14'
$ws.Range('E3').Value = '  -2.10%  '
# Row 4
$ws.Range('E4').Value = '  +0.12%  '
# Row 5
$ws.Range('D5').Value = '''601.38'
$ws.Range('E5').Value = '  -1.41%  '
# Row 6
$ws.Range('D6').Value = '''150.87'
$ws.Range('E6').Value = '  -5.34%  '
# Row 7
$ws.Range('E7').Value = '  +0.13%  '
# Row 8
$ws.Range('D8').Value = '3.164.82'
$ws.Range('E8').Value = '  -2.12%  '
# Row 9
$ws.Range('E9').Value = '  -3.26%  '
# Row 10
$ws.Range('E10').Value = '  -4.85%  '
# Row 11
$ws.Range('D11').Value = '''5.65'
$ws.Range('E11').Value = '  -1.61%  '
# Row 12
$ws.Range('D12').Value = '''0.479'
$ws.Range('E12').Value = '  -5.06%  '
# Row 13
$ws.Range('D13').Value = '''0.0000260'
$ws.Range('E13').Value = '  -4.18%  '
# Row 14
$ws.Range('D14').Value = '''37.14'
$ws.Range('E14').Value = '  -4.68%  '
# Row 15
$ws.Range('D15').Value = '3.668.74'
$ws.Range('E15').Value = '  -2.58%  '
# Row 16
$ws.Range('D16').Value = '64.760.53'
$ws.Range('E16').Value = '  -3.01%  '
# Row 17
$ws.Range('D17').Value = '3.191.62'
$ws.Range('E17').Value = '  -1.36%  '
# Row 18
$ws.Range('E18').Value = '  +0.53%  '
# Row 19
$ws.Range('D19').Value = '''7.05'
$ws.Range('E19').Value = '  -4.65%  '
# Row 20
$ws.Range('D20').Value = '''484.32'
$ws.Range('E20').Value = '  -5.14%  '
# Row 21
$ws.Range('D21').Value = '''14.88'
$ws.Range('E21').Value = '  -2.25%  '
# Row 22
$ws.Range('D22').Value = '''0.718'
$ws.Range('E22').Value = '  -2.54%  '
# Row 23
$ws.Range('D23').Value = '''7.81'
$ws.Range('E23').Value = '  -2.76%  '
# Row 24
$ws.Range('D24').Value = '''13.99'
$ws.Range('E24').Value = '  -4.13%  '
# Row 25
$ws.Range('D25').Value = '''85.17'
$ws.Range('E25').Value = '  +0.11%  '
# Row 26
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.13%  '
# Row 27
$ws.Range('E27').Value = '  -2.66%  '
# Row 28
$ws.Range('D28').Value = '''8.77'
$ws.Range('E28').Value = '  -4.42%  '
# Row 29
$ws.Range('E29').Value = '  -4.37%  '
# Row 30
$ws.Range('D30').Value = '''0.122'
$ws.Range('E30').Value = '  -0.51%  '
# Row 31
$ws.Range('D31').Value = '''7.13'
$ws.Range('E31').Value = '  +1.77%  '
# Row 32
$ws.Range('D32').Value = '''2.74'
$ws.Range('E32').Value = '  -6.65%  '
# Row 33
$ws.Range('E33').Value = '  -0.12%  '
# Row 34
$ws.Range('D34').Value = '''26.93'
$ws.Range('E34').Value = '  -4.53%  '
# Row 35
$ws.Range('E35').Value = '  -5.64%  '
# Row 36
$ws.Range('E36').Value = '  -5.44%  '
# Row 37
$ws.Range('D37').Value = '''54.87'
$ws.Range('E37').Value = '  -1.47%  '
# Row 38
$ws.Range('D38').Value = '''3.27'
$ws.Range('E38').Value = '  +5.88%  '
# Row 39
$ws.Range('D39').Value = '0.0₃0747'
$ws.Range('E39').Value = '  -3.41%  '
# Row 40
$ws.Range('D40').Value = '''459.99'
$ws.Range('E40').Value = '  -8.38%  '
# Row 41
$ws.Range('E41').Value = '  -3.04%  '
# Row 42
$ws.Range('D42').Value = '''0.0404'
$ws.Range('E42').Value = '  -4.19%  '
# Row 43
$ws.Range('E43').Value = '  -2.13%  '
# Row 44
$ws.Range('E44').Value = '  -0.34%  '
# Row 45
$ws.Range('D45').Value = '2.896.48'
$ws.Range('E45').Value = '  -0.79%  '
# Row 46
$ws.Range('E46').Value = '  -7.15%  '
# Row 47
$ws.Range('D47').Value = '''27.14'
$ws.Range('E47').Value = '  -3.68%  '
# Row 48
$ws.Range('E48').Value = '  -0.03%  '
# Row 49
$ws.Range('E49').Value = '  -2.98%  '
# Row 50
$ws.Range('E50').Value = '  -0.05%  '
# Row 51
$ws.Range('D51').Value = '''119.80'
$ws.Range('E51').Value = '  -2.19%  '
